# chore: update Sheets via scheduled runner
# Refresh computed market-board price/profit figures (columns H-N) on a
# number of rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 722.2222
$ws.Range("I4").Value = 722.2222
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 722.2222
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -608.2222

$ws.Range("H74").Value = 4571.4287
$ws.Range("I74").Value = 4600
$ws.Range("K74").Value = 4600
$ws.Range("M74").Value = -3664

$ws.Range("H77").Value = 4571.4287
$ws.Range("I77").Value = 4600
$ws.Range("K77").Value = 23000
$ws.Range("M77").Value = -18320

$ws.Range("H135").Value = 18520376
$ws.Range("I135").Value = 870.2381
$ws.Range("J135").Value = 83338650
$ws.Range("K135").Value = 7832.142900000001
$ws.Range("L135").Value = 750047850
$ws.Range("M135").Value = -5297.142900000001
$ws.Range("N135").Value = -750052920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4214.9375
$ws.Range("I74").Value = 731.5
$ws.Range("J74").Value = 6305
$ws.Range("K74").Value = 731.5
$ws.Range("L74").Value = 6305
$ws.Range("M74").Value = 142.5
$ws.Range("N74").Value = -8053

$ws.Range("H77").Value = 4214.9375
$ws.Range("I77").Value = 731.5
$ws.Range("J77").Value = 6305
$ws.Range("K77").Value = 3657.5
$ws.Range("L77").Value = 31525
$ws.Range("M77").Value = 710.5
$ws.Range("N77").Value = -40261

$ws.Range("H88").Value = 2854.4443
$ws.Range("I88").Value = 2998.5715
$ws.Range("J88").Value = 2350
$ws.Range("K88").Value = 2998.5715
$ws.Range("L88").Value = 2350
$ws.Range("M88").Value = -2592.5715
$ws.Range("N88").Value = -3162

$ws.Range("H91").Value = 2854.4443
$ws.Range("I91").Value = 2998.5715
$ws.Range("J91").Value = 2350
$ws.Range("K91").Value = 2998.5715
$ws.Range("L91").Value = 2350
$ws.Range("M91").Value = -1594.5715
$ws.Range("N91").Value = -5158

$ws.Range("H102").Value = 2700
$ws.Range("I102").Value = 2400
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2400
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -778
$ws.Range("N102").Value = -6244

$ws.Range("H132").Value = 1439781.9
$ws.Range("I132").Value = 1918431.5
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 5755294.5
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -5752764.5
$ws.Range("N132").Value = -16559

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").ClearContents()
$ws.Range("N61").Value = 0

$ws.Range("H107").Value = 232271.02
$ws.Range("I107").Value = 321388.97
$ws.Range("J107").Value = 1612.7646
$ws.Range("K107").Value = 321388.97
$ws.Range("L107").Value = 1612.7646
$ws.Range("M107").Value = -319468.97
$ws.Range("N107").Value = -5452.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1011.2308
$ws.Range("I35").Value = 1011.2308
$ws.Range("K35").Value = 1011.2308
$ws.Range("M35").Value = -717.2308

$ws.Range("H58").Value = 3667.2
$ws.Range("I58").Value = 1472.25
$ws.Range("J58").Value = 4465.364
$ws.Range("K58").Value = 1472.25
$ws.Range("L58").Value = 4465.364
$ws.Range("M58").Value = -1269.25
$ws.Range("N58").Value = -4871.364

$ws.Range("H134").Value = 2885.3667
$ws.Range("I134").Value = 2912.8928
$ws.Range("K134").Value = 8738.678400000001
$ws.Range("M134").Value = -6203.678400000001

$ws.Range("H136").Value = 3667.2
$ws.Range("I136").Value = 1472.25
$ws.Range("J136").Value = 4465.364
$ws.Range("K136").Value = 4416.75
$ws.Range("L136").Value = 13396.092
$ws.Range("M136").Value = -1866.75
$ws.Range("N136").Value = -18496.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 16129780
$ws.Range("I5").Value = 676.95
$ws.Range("J5").Value = 45455424
$ws.Range("K5").Value = 2030.85
$ws.Range("L5").Value = 136366272
$ws.Range("M5").Value = -1918.85
$ws.Range("N5").Value = -136366496

$ws.Range("H62").Value = 3664.6667
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 4997
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 14991
$ws.Range("M62").Value = -2314
$ws.Range("N62").Value = -16363

$ws.Range("H65").Value = 3664.6667
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 4997
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 44973
$ws.Range("M65").Value = -5568
$ws.Range("N65").Value = -51837

$ws.Range("H122").Value = 12756406
$ws.Range("J122").Value = 1489983.5
$ws.Range("L122").Value = 13409851.5
$ws.Range("N122").Value = -13414751.5

$ws.Range("H135").Value = 16129780
$ws.Range("I135").Value = 676.95
$ws.Range("J135").Value = 45455424
$ws.Range("K135").Value = 6092.55
$ws.Range("L135").Value = 409098816
$ws.Range("M135").Value = -3557.55
$ws.Range("N135").Value = -409103886

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1544.9333
$ws.Range("I126").Value = 1272.8334
$ws.Range("J126").Value = 2633.3333
$ws.Range("K126").Value = 3818.5002
$ws.Range("L126").Value = 7899.999899999999
$ws.Range("M126").Value = -1348.5002
$ws.Range("N126").Value = -12839.9999

$ws.Range("H132").Value = 2588.1562
$ws.Range("I132").Value = 1954.1578
$ws.Range("K132").Value = 5862.4734
$ws.Range("M132").Value = -3332.4734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 395
$ws.Range("I22").Value = 366.66666
$ws.Range("J22").Value = 412
$ws.Range("K22").Value = 366.66666
$ws.Range("L22").Value = 412
$ws.Range("M22").Value = -71.66665999999998
$ws.Range("N22").Value = -1002

$ws.Range("H27").Value = 395
$ws.Range("I27").Value = 366.66666
$ws.Range("J27").Value = 412
$ws.Range("K27").Value = 366.66666
$ws.Range("L27").Value = 412
$ws.Range("M27").Value = -259.66666
$ws.Range("N27").Value = -626

$ws.Range("H55").Value = 172.82608
$ws.Range("I55").Value = 145.9375
$ws.Range("J55").Value = 234.28572
$ws.Range("K55").Value = 145.9375
$ws.Range("L55").Value = 234.28572
$ws.Range("M55").Value = 27.0625
$ws.Range("N55").Value = -580.28572

$ws.Range("H122").Value = 5053.3022
$ws.Range("I122").Value = 5617.6177
$ws.Range("K122").Value = 16852.8531
$ws.Range("M122").Value = -14402.8531

$ws.Range("H132").Value = 6218.3706
$ws.Range("I132").Value = 7484.3687
$ws.Range("J132").Value = 3211.625
$ws.Range("K132").Value = 22453.1061
$ws.Range("L132").Value = 9634.875
$ws.Range("M132").Value = -19923.1061
$ws.Range("N132").Value = -14694.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H136").Value = 1259.1428
$ws.Range("I136").Value = 1119
$ws.Range("J136").Value = 2100
$ws.Range("K136").Value = 3357
$ws.Range("L136").Value = 6300
$ws.Range("M136").Value = -807
$ws.Range("N136").Value = -11400
